$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5's phone number was stored as text; normalize it to a real number.
$ws.Range("A5").Value = 79174463

# Add new payment row 6: 79174460 (Cash) 2025-08-20T09:41:48
# Column A keeps phone numbers stored as text in this sheet's source data,
# so force a text value for A6 (the numeric-looking "@" format is only used
# transiently to stop Excel from auto-coercing the value to a number, then
# cleared so the cell keeps the sheet's default, unstyled formatting).
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "79174460"
$ws.Range("A6").ClearFormats()

$ws.Range("B6").Value = 50
$ws.Range("C6").Value = 7.5
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 42.5
$ws.Range("G6").Value = "Cash"
$ws.Range("H6").Value = "2025-08-20T09:41:48"
